# refactor currency conversion, now explicit source and target amounts
#
# currency_conversions sheet used to have a single "foreign_amount" column;
# split this into explicit "source_amount" / "target_amount" (+ a new
# "target_fees" column alongside the existing "source_fees") so both legs of
# a conversion are recorded independently.
#
# Old layout: date | foreign_amount | source_fees | source_currency | target_currency | comment
# New layout: date | source_amount  | source_fees | source_currency | target_amount | target_fees | target_currency | comment

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("currency_conversions")

# B1 "foreign_amount" becomes "source_amount" (same column, renamed).
$ws.Range("B1").Value = "source_amount"

# Make room for the two new columns ("target_amount", "target_fees") right
# before the existing "target_currency" column, shifting it (and "comment")
# two columns to the right.
$ws.Range("E1").EntireColumn.Insert()
$ws.Range("E1").EntireColumn.Insert()

$ws.Range("E1").Value = "target_amount"
$ws.Range("F1").Value = "target_fees"

# This sheet becomes the active / selected tab of the workbook.
$ws.Activate()
